# Append 5 new data rows (265-269) to Sheet1, mirroring the existing
# "date" / "value" log rows already present in A2:B264.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("04/27/2021 18:41:16", 0.3333),
    @("04/27/2021 18:43:06", 0.339),
    @("04/27/2021 18:51:30", 0.3333),
    @("04/27/2021 18:54:44", 0.1667),
    @("04/27/2021 18:57:56", 0.3333)
)

$startRow = 265
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}
